$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 193.27777
$ws.Range("I5").Value = 52.916668
$ws.Range("J5").Value = 474
$ws.Range("K5").Value = 52.916668
$ws.Range("L5").Value = 474
$ws.Range("M5").Value = 62.083332
$ws.Range("N5").Value = -704

$ws.Range("H51").Value = 29621.432
$ws.Range("I51").Value = 9624.9375
$ws.Range("J51").Value = 44856.855
$ws.Range("K51").Value = 9624.9375
$ws.Range("L51").Value = 44856.855
$ws.Range("M51").Value = -9140.9375
$ws.Range("N51").Value = -45824.855

$ws.Range("H62").Value = 12124479
$ws.Range("I62").Value = 19051068
$ws.Range("J62").Value = 2949.75
$ws.Range("K62").Value = 19051068
$ws.Range("L62").Value = 2949.75
$ws.Range("M62").Value = -19050444
$ws.Range("N62").Value = -4197.75

$ws.Range("H65").Value = 12124479
$ws.Range("I65").Value = 19051068
$ws.Range("J65").Value = 2949.75
$ws.Range("K65").Value = 95255340
$ws.Range("L65").Value = 14748.75
$ws.Range("M65").Value = -95252220
$ws.Range("N65").Value = -20988.75

$ws.Range("H70").Value = 1227608.4
$ws.Range("J70").Value = 2117232.8
$ws.Range("L70").Value = 6351698.399999999
$ws.Range("N70").Value = -6352238.399999999

$ws.Range("H73").Value = 1227608.4
$ws.Range("J73").Value = 2117232.8
$ws.Range("L73").Value = 6351698.399999999
$ws.Range("N73").Value = -6353570.399999999

$ws.Range("H74").Value = 7399.5386
$ws.Range("I74").Value = 4839.2
$ws.Range("K74").Value = 4839.2
$ws.Range("M74").Value = -3903.2

$ws.Range("H77").Value = 7399.5386
$ws.Range("I77").Value = 4839.2
$ws.Range("K77").Value = 24196
$ws.Range("M77").Value = -19516

$ws.Range("H98").Value = 1731.55
$ws.Range("I98").Value = 1695.2354
$ws.Range("K98").Value = 1695.2354
$ws.Range("M98").Value = -197.2354

$ws.Range("H122").Value = 1731.55
$ws.Range("I122").Value = 1695.2354
$ws.Range("K122").Value = 5085.706200000001
$ws.Range("M122").Value = -2635.706200000001

$ws.Range("H129").Value = 2625.2307
$ws.Range("J129").Value = 2741.182
$ws.Range("L129").Value = 8223.545999999998
$ws.Range("N129").Value = -18223.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30086.906
$ws.Range("I32").Value = 30583.025
$ws.Range("K32").Value = 30583.025
$ws.Range("M32").Value = -30296.025

$ws.Range("H122").Value = 6379.6
$ws.Range("I122").Value = 6379.6
$ws.Range("K122").Value = 19138.8
$ws.Range("M122").Value = -16688.8

$ws.Range("H131").Value = 84205.57000000001
$ws.Range("J131").Value = 84205.57000000001
$ws.Range("L131").Value = 84205.57000000001
$ws.Range("N131").Value = -94285.57000000001

$ws.Range("H132").Value = 4475.278
$ws.Range("I132").Value = 3154.9
$ws.Range("K132").Value = 9464.700000000001
$ws.Range("M132").Value = -6934.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 52000
$ws.Range("J21").Value = 52000
$ws.Range("L21").Value = 52000
$ws.Range("N21").Value = -52472

$ws.Range("H134").Value = 5632.5
$ws.Range("I134").Value = 3807.1333
$ws.Range("K134").Value = 11421.3999
$ws.Range("M134").Value = -8886.3999

$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1977.5555
$ws.Range("I16").Value = 1974.75
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1974.75
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1687.75
$ws.Range("N16").Value = -2574

$ws.Range("H62").Value = 13724.5
$ws.Range("I62").Value = 4966
$ws.Range("K62").Value = 4966
$ws.Range("M62").Value = -4342

$ws.Range("H65").Value = 13724.5
$ws.Range("I65").Value = 4966
$ws.Range("K65").Value = 24830
$ws.Range("M65").Value = -21710

$ws.Range("H113").Value = 1977.5555
$ws.Range("I113").Value = 1974.75
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1974.75
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 195.25
$ws.Range("N113").Value = -6340

$ws.Range("H132").Value = 70550
$ws.Range("I132").Value = 8504.5
$ws.Range("K132").Value = 25513.5
$ws.Range("M132").Value = -22983.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 2016.6666
$ws.Range("I26").Value = 80
$ws.Range("J26").Value = 11700
$ws.Range("K26").Value = 240
$ws.Range("L26").Value = 35100
$ws.Range("M26").Value = 48
$ws.Range("N26").Value = -35676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5150
$ws.Range("I43").Value = 5150
$ws.Range("K43").Value = 5150
$ws.Range("M43").Value = -4999

$ws.Range("H122").Value = 1175.7273
$ws.Range("I122").Value = 1099.6666
$ws.Range("J122").Value = 1204.25
$ws.Range("K122").Value = 3298.9998
$ws.Range("L122").Value = 3612.75
$ws.Range("M122").Value = -848.9998000000001
$ws.Range("N122").Value = -8512.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3936.0588
$ws.Range("I22").Value = 2785.5386
$ws.Range("K22").Value = 2785.5386
$ws.Range("M22").Value = -2490.5386

$ws.Range("H27").Value = 3936.0588
$ws.Range("I27").Value = 2785.5386
$ws.Range("K27").Value = 2785.5386
$ws.Range("M27").Value = -2678.5386

$ws.Range("H40").Value = 10268.759
$ws.Range("I40").Value = 8946.333000000001
$ws.Range("K40").Value = 8946.333000000001
$ws.Range("M40").Value = -8810.333000000001

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("N51").Value = 0
$ws.Range("L51").ClearContents()

$ws.Range("H61").Value = 4300.25
$ws.Range("I61").Value = 3849.5
$ws.Range("K61").Value = 3849.5
$ws.Range("M61").Value = -3647.5

$ws.Range("H68").Value = 2762.4375
$ws.Range("J68").Value = 3420
$ws.Range("L68").Value = 3420
$ws.Range("N68").Value = -4918

$ws.Range("H71").Value = 2762.4375
$ws.Range("J71").Value = 3420
$ws.Range("L71").Value = 17100
$ws.Range("N71").Value = -24588

$ws.Range("H82").Value = 1130.4445
$ws.Range("J82").Value = 1008.1667
$ws.Range("L82").Value = 1008.1667
$ws.Range("N82").Value = -1730.1667

$ws.Range("H85").Value = 1130.4445
$ws.Range("J85").Value = 1008.1667
$ws.Range("L85").Value = 1008.1667
$ws.Range("N85").Value = -3504.1667

$ws.Range("H113").Value = 4300.25
$ws.Range("I113").Value = 3849.5
$ws.Range("K113").Value = 3849.5
$ws.Range("M113").Value = -1679.5

$ws.Range("H136").Value = 4731.5713
$ws.Range("I136").Value = 4427.212
$ws.Range("K136").Value = 13281.636
$ws.Range("M136").Value = -10731.636

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()

$ws.Range("H62").Value = 5475.75
$ws.Range("I62").Value = 5900
$ws.Range("J62").Value = 5334.3335
$ws.Range("K62").Value = 5900
$ws.Range("L62").Value = 5334.3335
$ws.Range("M62").Value = -5276
$ws.Range("N62").Value = -6582.3335

$ws.Range("H65").Value = 5475.75
$ws.Range("I65").Value = 5900
$ws.Range("J65").Value = 5334.3335
$ws.Range("K65").Value = 29500
$ws.Range("L65").Value = 26671.6675
$ws.Range("M65").Value = -26380
$ws.Range("N65").Value = -32911.6675

$ws.Range("H113").Value = 841.2414
$ws.Range("I113").Value = 964.2381
$ws.Range("J113").Value = 518.375
$ws.Range("K113").Value = 2892.7143
$ws.Range("L113").Value = 1555.125
$ws.Range("M113").Value = -722.7143000000001
$ws.Range("N113").Value = -5895.125

$ws.Range("H132").Value = 9083.134
$ws.Range("I132").Value = 7519.7
$ws.Range("K132").Value = 22559.1
$ws.Range("M132").Value = -20029.1

$ws.Range("H136").Value = 5922.846
$ws.Range("I136").Value = 2267.6
$ws.Range("J136").Value = 8207.375
$ws.Range("K136").Value = 6802.799999999999
$ws.Range("L136").Value = 24622.125
$ws.Range("M136").Value = -4252.799999999999
$ws.Range("N136").Value = -29722.125
